$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = -0.63
$ws.Range("K2").Value = -0.72
$ws.Range("L2").Value = -0.47
$ws.Range("J3").Value = 0.13
$ws.Range("K3").Value = -0.32
$ws.Range("L3").Value = 0.9
$ws.Range("H4").Value = -0
$ws.Range("H5").Value = -0
$ws.Range("K6").Value = -0.03
$ws.Range("L6").Value = 0.6899999999999999
$ws.Range("J7").Value = 0.13
$ws.Range("K7").Value = -1.52
$ws.Range("L7").Value = 1.71
$ws.Range("H8").Value = -0
$ws.Range("H9").Value = -0
$ws.Range("I9").Value = -0
$ws.Range("H10").Value = -0
$ws.Range("H11").Value = -0
$ws.Range("H12").Value = -0
$ws.Range("I12").Value = -0
$ws.Range("H13").Value = -0
$ws.Range("I13").Value = -0
$ws.Range("J14").Value = -5.95
$ws.Range("K14").Value = -7.45
$ws.Range("L14").Value = -4.69
$ws.Range("J15").Value = -4.13
$ws.Range("K15").Value = -5.9
$ws.Range("L15").Value = -2.55
$ws.Range("J16").Value = -6.57
$ws.Range("K16").Value = -8.66
$ws.Range("L16").Value = -4.48
$ws.Range("J17").Value = -6.84
$ws.Range("K17").Value = -8.98
$ws.Range("L17").Value = -5.16
$ws.Range("J18").Value = -6.87
$ws.Range("K18").Value = -8.390000000000001
$ws.Range("L18").Value = -5.35
$ws.Range("J19").Value = -6.46
$ws.Range("K19").Value = -7.34
$ws.Range("L19").Value = -5.01
$ws.Range("J20").Value = -9.640000000000001
$ws.Range("K20").Value = -10.33
$ws.Range("L20").Value = -8.960000000000001
$ws.Range("J21").Value = -7.43
$ws.Range("K21").Value = -8.42
$ws.Range("L21").Value = -5.94
$ws.Range("J22").Value = -6.91
$ws.Range("K22").Value = -8.99
$ws.Range("L22").Value = -4.83
$ws.Range("J23").Value = -6.94
$ws.Range("K23").Value = -9.33
$ws.Range("L23").Value = -3.91
$ws.Range("J24").Value = -5.58
$ws.Range("K24").Value = -5.58
$ws.Range("L24").Value = -5.58
$ws.Range("J25").Value = -7.42
$ws.Range("K25").Value = -8.41
$ws.Range("L25").Value = -6.28
$ws.Range("J26").Value = -0.54
$ws.Range("K26").Value = -0.97
$ws.Range("L26").Value = -0.19
$ws.Range("J27").Value = -2.59
$ws.Range("L27").Value = -2.13
$ws.Range("J28").Value = -3.28
$ws.Range("K28").Value = -4.12
$ws.Range("L28").Value = -2.44
$ws.Range("J29").Value = -1.83
$ws.Range("K29").Value = -2.62
$ws.Range("L29").Value = -0.6
$ws.Range("J30").Value = -0.99
$ws.Range("K30").Value = -1.53
$ws.Range("L30").Value = -0.45
$ws.Range("J31").Value = -3.38
$ws.Range("K31").Value = -4.28
$ws.Range("L31").Value = -2.06
$ws.Range("J32").Value = -3.1
$ws.Range("K32").Value = -3.77
$ws.Range("L32").Value = -2.44
$ws.Range("J33").Value = -1.79
$ws.Range("K33").Value = -3.11
$ws.Range("L33").Value = -0.16
$ws.Range("J34").Value = -2.5
$ws.Range("K34").Value = -3.42
$ws.Range("L34").Value = -1.58
$ws.Range("J35").Value = -1.93
$ws.Range("K35").Value = -3.57
$ws.Range("L35").Value = -1.1
$ws.Range("J36").Value = -2.43
$ws.Range("K36").Value = -2.43
$ws.Range("L36").Value = -2.43
$ws.Range("J37").Value = -2.12
$ws.Range("K37").Value = -2.81
$ws.Range("L37").Value = -1.62
$ws.Range("K38").Value = -0.73
$ws.Range("L38").Value = 1.08
$ws.Range("J39").Value = -0.55
$ws.Range("K39").Value = -0.89
$ws.Range("J40").Value = -0.6899999999999999
$ws.Range("K40").Value = -2.02
$ws.Range("L40").Value = 0.64
$ws.Range("J41").Value = -0.61
$ws.Range("K41").Value = -2.78
$ws.Range("L41").Value = 0.55
$ws.Range("J42").Value = 0.98
$ws.Range("K42").Value = 0.53
$ws.Range("L42").Value = 1.43
$ws.Range("J43").Value = 0.01
$ws.Range("K43").Value = -0.78
$ws.Range("L43").Value = 1.55
$ws.Range("J44").Value = -0.01
$ws.Range("K44").Value = -0.48
$ws.Range("L44").Value = 0.45
$ws.Range("J45").Value = -1
$ws.Range("K45").Value = -1.85
$ws.Range("L45").Value = -0.18
$ws.Range("J46").Value = -0.37
$ws.Range("L46").Value = 0.22
$ws.Range("K47").Value = -1.15
$ws.Range("L47").Value = 0.63
$ws.Range("J48").Value = -0.8100000000000001
$ws.Range("K48").Value = -0.8100000000000001
$ws.Range("L48").Value = -0.8100000000000001
$ws.Range("K49").Value = -0.65
$ws.Range("J50").Value = -0.02
$ws.Range("K50").Value = -0.77
$ws.Range("L50").Value = 0.79
$ws.Range("J51").Value = -0.9
$ws.Range("K51").Value = -1.98
$ws.Range("L51").Value = -0.34
$ws.Range("J52").Value = -1.33
$ws.Range("K52").Value = -2.28
$ws.Range("L52").Value = -0.38
$ws.Range("J53").Value = -1.44
$ws.Range("K53").Value = -1.86
$ws.Range("L53").Value = -1.13
$ws.Range("J54").Value = -0.07000000000000001
$ws.Range("K54").Value = -1.21
$ws.Range("L54").Value = 1.07
$ws.Range("J55").Value = -1
$ws.Range("K55").Value = -1.86
$ws.Range("L55").Value = -0.3
$ws.Range("J56").Value = -1.44
$ws.Range("K56").Value = -1.87
$ws.Range("L56").Value = -1.02
$ws.Range("J57").Value = -1.41
$ws.Range("K57").Value = -2.5
$ws.Range("L57").Value = -0.72
$ws.Range("J58").Value = -0.22
$ws.Range("K58").Value = -0.37
$ws.Range("L58").Value = -0.06
$ws.Range("J59").Value = -0.97
$ws.Range("K59").Value = -1.46
$ws.Range("L59").Value = -0.09
$ws.Range("J60").Value = -1.71
$ws.Range("K60").Value = -1.71
$ws.Range("L60").Value = -1.71
$ws.Range("J61").Value = -1.07
$ws.Range("K61").Value = -1.63
$ws.Range("L61").Value = -0.5600000000000001
